# Add a new "Save" column (H) to the s_vals sheet, mirroring the existing
# header style used by the other stat columns (e.g. G1 "sum").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1: bold, centered,
# top-aligned, thin border) onto the new header cell H1 so it reuses the
# same cell style rather than minting a new one.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header text and the data value for row 2.
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
